$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.713.58"
$ws.Range("E2").Value = "  -3.00%  "

$ws.Range("D3").Value = "1.773.06"
$ws.Range("E3").Value = "  -3.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.51"
$ws.Range("E5").Value = "  -6.87%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5003"
$ws.Range("E7").Value = "  -4.86%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2479"
$ws.Range("E8").Value = "  -22.12%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06166"
$ws.Range("E9").Value = "  -9.19%  "

$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.800.92"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06692"
$ws.Range("E11").Value = "  -13.58%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.62"
$ws.Range("E12").Value = "  -21.98%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6075"
$ws.Range("E13").Value = "  -22.48%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "78.53"
$ws.Range("E14").Value = "  -10.59%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.356"
$ws.Range("E15").Value = "  -13.05%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "25.755.75"
$ws.Range("E18").Value = "  -2.94%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("E19").Value = "  -18.62%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.023.76"
$ws.Range("E20").Value = "  -2.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006314"
$ws.Range("E21").Value = "  -20.59%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.900"
$ws.Range("E22").Value = "  -15.53%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.175"
$ws.Range("E23").Value = "  -13.38%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.014"
$ws.Range("E24").Value = "  -13.97%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "131.44"
$ws.Range("E25").Value = "  -7.46%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.886"
$ws.Range("E26").Value = "  -14.51%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.45"
$ws.Range("E27").Value = "  -14.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.372"
$ws.Range("E28").Value = "  -18.57%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "98.64"
$ws.Range("E29").Value = "  -11.54%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08213"
$ws.Range("E30").Value = "  -5.71%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.592"
$ws.Range("E31").Value = "  -13.66%  "

$ws.Range("B32").Value = "Frax"
$ws.Range("C32").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.002"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.724"
$ws.Range("E33").Value = "  -4.73%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04288"
$ws.Range("E34").Value = "  -12.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.139"
$ws.Range("E35").Value = "  -22.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.030"
$ws.Range("E36").Value = "  -9.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6202"
$ws.Range("E37").Value = "  -14.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.774"
$ws.Range("E38").Value = "  -10.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.111"
$ws.Range("E39").Value = "  -5.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.25"
$ws.Range("E41").Value = "  -7.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01453"
$ws.Range("E42").Value = "  -17.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7819"
$ws.Range("E43").Value = "  -12.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3867"
$ws.Range("E44").Value = "  -19.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.159"
$ws.Range("E45").Value = "  -12.84%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05231"
$ws.Range("E46").Value = "  -10.52%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.149"
$ws.Range("E47").Value = "  -19.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.48"
$ws.Range("E48").Value = "  -12.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.494"
$ws.Range("E50").Value = "  -16.19%  "

$ws.Range("B51").Value = "TrueUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.03%  "
